$wb = $excel.ActiveWorkbook

# --- Sheet "Orders": update F21 and append rows 22..41 -----------------
$ws = $wb.Worksheets.Item("Orders")

function Set-TextValue($cell, $value) {
    if ($null -eq $value) { return }
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# F21: "1" -> "10" (kept as text, matching the rest of the column)
Set-TextValue $ws.Cells.Item(21, 6) "10"

$newRows = @(
    @{ Row = 22; A = '6'; C = '633_干花安娜深红_undefined_undefined_1stem'; F = '15' }
    @{ Row = 23; A = '1'; C = '147_娜欧米_Red Naomi_Rosa rugosa Thunb._20stems'; F = '6' }
    @{ Row = 24; A = $null; C = '192_粉荔枝_Pink Ohara_Rosa rugosa Thunb._20stems'; F = '12' }
    @{ Row = 25; A = $null; C = '148_坦尼克_Tineke_Rosa rugosa Thunb._20stems'; F = '10' }
    @{ Row = 26; A = $null; C = '197_粉红雪山_Sweet Avalanche_Rosa rugosa Thunb._20stems'; F = '14' }
    @{ Row = 27; A = $null; C = '12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g'; F = '10' }
    @{ Row = 28; A = '2'; C = '197_粉红雪山_Sweet Avalanche_Rosa rugosa Thunb._20stems'; F = '6' }
    @{ Row = 29; A = $null; C = '274_仙子之吻_undefined_Rosa rugosa Thunb._10stems'; F = '16' }
    @{ Row = 30; A = $null; C = '12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g'; F = '5' }
    @{ Row = 31; A = $null; C = '3_波浪白洋桔梗_Wavy White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g'; F = '10' }
    @{ Row = 32; A = $null; C = '625_多丁紫蝴蝶_undefined_undefined_1bunch'; F = '5' }
    @{ Row = 33; A = $null; C = '412_紫罗兰粉_violet pink_undefined_1bunch'; F = '10' }
    @{ Row = 34; A = $null; C = '512_松虫草粉_scabiosa pink_undefined_1bunch'; F = '5' }
    @{ Row = 35; A = $null; C = '419_松虫草红_scabiosa watermelon_undefined_1bunch'; F = '7' }
    @{ Row = 36; A = $null; C = '100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem'; F = '30' }
    @{ Row = 37; A = $null; C = '107_绣球单瓣浅粉_Hydrangea Light Pink S_Hydrangea L._1stem'; F = '10' }
    @{ Row = 38; A = $null; C = '302_彩星 浅粉_Tinted Gypso light pink_undefined_0.5kg'; F = '10' }
    @{ Row = 39; A = '3'; C = '107_绣球单瓣浅粉_Hydrangea Light Pink S_Hydrangea L._1stem'; F = '35' }
    @{ Row = 40; A = $null; C = '571_大飞燕浅紫_undefined_undefined_1bunch'; F = '25' }
    @{ Row = 41; A = $null; C = '647_海棠果红_undefined_undefined_1bunch'; F = $null }
)

foreach ($r in $newRows) {
    Set-TextValue $ws.Cells.Item($r.Row, 1) $r.A   # PackageID
    Set-TextValue $ws.Cells.Item($r.Row, 3) $r.C   # FlowerName
    Set-TextValue $ws.Cells.Item($r.Row, 6) $r.F   # Number
}

# --- Sheet "Summary": update the concatenated PackageID code in G2 -----
$ws2 = $wb.Worksheets.Item("Summary")
Set-TextValue $ws2.Cells.Item(2, 7) "02424158281012115205540502050501010156121014106165105105730101035250"

Write-Host "Edit applied: Orders!F21:F41 updated, Summary!G2 updated."
